$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the buurtcode / gemeente values in row 2
$ws.Range("A2").Value = "BU15810004"
$ws.Range("B2").Value = "unknown"

# Bold, small, dark-grey Arial label style on the buurtcode cell
$font = $ws.Range("A2").Font
$font.Name = "Arial"
$font.Bold = $true
$font.Size = 7
$font.Color = 2171169

# Move the active selection to B2
[void]$ws.Range("B2").Select()

# Portrait page orientation for printing
$ws.PageSetup.Orientation = 1
